# Update today work hours
# -------------------------------------------------------------------
# Reproduces, via Excel COM automation, the weekly-timesheet update:
#  1) the "planned hours" cells (rows 3:4) and "actual hours" cells
#     (rows 6:7) swap their pre-existing highlight styles so that
#     "planned" consistently uses the yellow/Note style and "actual"
#     consistently uses the green/Good style;
#  2) "today" (rows 6:7, the Friday columns K/L/M) gets its actual
#     worked hours filled in;
#  3) a whole new (blank) week -- rows 10-15, mirroring the first
#     week's 2-7 layout -- is appended below.
# -------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-CellFormat($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy()
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------
# 1) Swap the "planned" (rows 3:4) / "actual" (rows 6:7) highlight
#    styles. Stash the "actual" style (from D6) in a scratch cell
#    first, since D3/D6 get overwritten below.
# ---------------------------------------------------------------
Copy-CellFormat "D6" "O1"
Copy-CellFormat "D3" "O2"

foreach ($addr in @("D3", "F3", "H3", "J3", "L3", "D4", "F4", "H4", "J4", "L4")) {
    Copy-CellFormat "O1" $addr
}

foreach ($addr in @("D6", "D7")) {
    Copy-CellFormat "O2" $addr
}

$ws.Range("O1:O2").Clear()

# ---------------------------------------------------------------
# 2) Fill in today's actual worked hours (Friday column, rows 6:7).
# ---------------------------------------------------------------
$ws.Range("K6").Value = 0.64583333333333337
$ws.Range("L6").Value = 4.25
$ws.Range("M6").Value = 23
$ws.Range("K7").Value = 0.82291666666666663

# ---------------------------------------------------------------
# 3) Append next week's (blank) timesheet block, rows 10-15,
#    mirroring rows 2-7 (rows 9/13 stay blank spacer rows).
#    Merge the vertically-spanning cells *before* copying formats
#    into them, so the copied format lands cleanly on the already
#    merged range instead of Excel minting a fresh per-cell style.
# ---------------------------------------------------------------
$ws.Range("B11:B12").Merge()
$ws.Range("D11:D12").Merge()
$ws.Range("F11:F12").Merge()
$ws.Range("H11:H12").Merge()
$ws.Range("J11:J12").Merge()
$ws.Range("L11:L12").Merge()
$ws.Range("M11:M12").Merge()

$ws.Range("B14:B15").Merge()
$ws.Range("D14:D15").Merge()
$ws.Range("F14:F15").Merge()
$ws.Range("H14:H15").Merge()
$ws.Range("J14:J15").Merge()
$ws.Range("L14:L15").Merge()
$ws.Range("M14:M15").Merge()

# Row 10 - date header (mirrors row 2, with next week's dates)
foreach ($col in @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M")) {
    Copy-CellFormat "$($col)2" "$($col)10"
}
$ws.Range("C10").Value = 43444
$ws.Range("E10").Value = 43445
$ws.Range("G10").Value = 43446
$ws.Range("I10").Value = 43447
$ws.Range("K10").Value = 43448
$ws.Range("M10").Value = "сумма"

# Rows 11-12 - "planned hours" block (mirrors rows 3-4, left blank)
foreach ($col in @("B", "C", "D", "E", "F", "H", "I", "J", "K", "L", "M")) {
    Copy-CellFormat "$($col)3" "$($col)11"
    Copy-CellFormat "$($col)4" "$($col)12"
}
$ws.Range("B11").Value = "Планируемые часы работы"

# Rows 14-15 - "actual hours" block (mirrors rows 6-7, left blank)
foreach ($col in @("B", "C", "D", "E", "F", "H", "I", "J", "K", "L", "M")) {
    Copy-CellFormat "$($col)6" "$($col)14"
    Copy-CellFormat "$($col)7" "$($col)15"
}
$ws.Range("B14").Value = "Фактические часы работы"

# ---------------------------------------------------------------
# 4) Update the view: selection moves to L13.
# ---------------------------------------------------------------
$ws.Range("L13").Select()
